$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing price cells (harga beli revisions) ---
$ws.Range("B13").Value = 20000
$ws.Range("B47").Value = 90000
$ws.Range("B62").Value = 140000
$ws.Range("B66").Value = 55000

# --- Rows 69:71 were blank spacer rows. Insert 3 fresh rows in their
#     place (this clones the cell formatting of row 68 above, same as a
#     manual Excel "Insert Row" would do), then drop the old blank
#     spacer rows which got pushed down to 72:74 so the row count /
#     numbering below stays exactly as it was. ---
$ws.Range("69:71").Insert()
$ws.Range("72:74").Delete()
$ws.Rows.Item(69).RowHeight = 15.75
$ws.Rows.Item(70).RowHeight = 15.75
$ws.Rows.Item(71).RowHeight = 15.75

# --- Row 69: Tensimeter digital - Sinocare BSX 516 ---
$ws.Range("E69").Value = "Sinocare BSX 516"
$ws.Range("F69").Value = "Tensimeter digital"
$ws.Range("A69").Formula = '=F69 & " - " & E69'
$ws.Range("B69").Value = 200000
$ws.Range("C69").Formula = '=UPPER(LEFT(F69,4) & "-" & LEFT(E69,2) & "-" & TEXT(ROW(A69)-1,"000"))'
$ws.Range("D69").Value = 10

# --- Row 70: Saturasi Oksigen - Omicron ---
$ws.Range("E70").Value = "Omicron"
$ws.Range("F70").Value = "Saturasi Oksigen"
$ws.Range("A70").Formula = '=F70 & " - " & E70'
$ws.Range("B70").Value = 120000
$ws.Range("C70").Formula = '=UPPER(LEFT(F70,4) & "-" & LEFT(E70,2) & "-" & TEXT(ROW(A70)-1,"000"))'
$ws.Range("D70").Value = 10

# --- Row 71: Saturasi Oksigen - GEA ---
$ws.Range("E71").Value = "GEA"
$ws.Range("F71").Value = "Saturasi Oksigen"
$ws.Range("A71").Formula = '=F71 & " - " & E71'
$ws.Range("B71").Value = 210000
$ws.Range("C71").Formula = '=UPPER(LEFT(F71,4) & "-" & LEFT(E71,2) & "-" & TEXT(ROW(A71)-1,"000"))'
$ws.Range("D71").Value = 2
